# Replace the single-cell "questions" payload with a pretty-printed JSON
# rendering (it used to be a Python dict/list repr()), moving it from A2
# up to A1, and drop the old placeholder numeric cell in A1 + its
# bold/bordered/centered style so the sheet ends up with a single, plainly
# styled A1 cell holding the new text.

$b64 = "cXVlc3Rpb25zID0gWwogICAgewogICAgICAgICJ0aXRsZSI6ICJZb3VyIGNsb3VkIGhvc3RpbmcgY2hhcmdlcyB5b3UgYSBsb3Qgb2YgbW9uZXkgZm9yIG5ldHdvcmsgdHJhZmZpYy4gWW91IHRyb3VibGVzaG9vdCB0aGUgaXNzdWUgYW5kIG5vdGljZSB0aGF0IHZlcnkgbGFyZ2UgSmF2YXNjcmlwdCBmaWxlcyB0YWtlIHVwIG1vc3Qgb2YgdGhlIGJhbmR3aWR0aC5XaGF0IHN0cmF0ZWdpZXMgc2hvdWxkIHlvdSBjb25zaWRlciBmb3IgdXNpbmcgbGVzcyBiYW5kd2lkdGg/IiwKICAgICAgICAicXVlc190eXBlIjogMTUsCiAgICAgICAgIm9wdGlvbnMiOiBbCiAgICAgICAgICAgICJTdHJpcCBjb21tZW50cyBhbmQgZXhjZXNzIHdoaXRlc3BhY2UgYmVmb3JlIGRlcGxveWluZy4iLAogICAgICAgICAgICAiRGVwbG95IHRoZSBwYWNrYWdlIHdpdGggbnBtIHJ1biBkZXBsb3kuIiwKICAgICAgICAgICAgIlNlcnZlIHRoZSBmaWxlcyBvdmVyIEFwYWNoZSBpbnN0ZWFkIG9mIE5naW54LiIsCiAgICAgICAgICAgICJDb21wcmVzcyB0aGUgZmlsZSBpbiB0cmFuc2l0IHdpdGggZ3ppcC4iLAogICAgICAgICAgICAiU2NhbiB0aGUgY29kZSB3aXRoIGEgbGludGVyLiIsCiAgICAgICAgICAgICJJbXBsZW1lbnQgdHlwZSBwcm90ZWN0aW9ucyB3aXRoIHRoZSBUeXBlc2NyaXB0IGNvbXBpbGVyLiIKICAgICAgICBdLAogICAgICAgICJzY29yZSI6IFsKICAgICAgICAgICAgIlN0cmlwIGNvbW1lbnRzIGFuZCBleGNlc3Mgd2hpdGVzcGFjZSBiZWZvcmUgZGVwbG95aW5nLiIsCiAgICAgICAgICAgICJDb21wcmVzcyB0aGUgZmlsZSBpbiB0cmFuc2l0IHdpdGggZ3ppcC4iCiAgICAgICAgXQogICAgfSwKICAgIHsKICAgICAgICAidGl0bGUiOiAiWW91IGFyZSByZXZpZXdpbmcgYSBjb2xsZWFndWVcdTIwMTlzIGNvZGUgYW5kIG5vdGljZSB0aGF0IHRoZXkgaGF2ZSB3cml0dGVuIGEgbG9vcCBzdHJ1Y3R1cmVkIGFzIHNob3duIGJlbG93LkhvdyBzaG91bGQgeW91IHJld3JpdGUgdGhpcyBtb3JlIGNvbmNpc2VseSB3aGlsZSBrZWVwaW5nIHRoZSBzYW1lIGZ1bmN0aW9uYWxpdHk/IyB3YWl0IGZvciBuZXR3b3JrIHRvIGNvbWUgb25saW5lXG53aGlsZSBUcnVlOlxuICAgIGlmIG5ldHdvcmsuaXNfb25saW5lKCk6XG4gICAgICAgIGJyZWFrXG4gICAgc2xlZXAoMSkgIyB3YWl0IG9uZSBzZWNvbmQsIHRoZW4gdHJ5IGFnYWluLiIsCiAgICAgICAgInF1ZXNfdHlwZSI6IDIsCiAgICAgICAgIm9wdGlvbnMiOiBbCiAgICAgICAgICAgICJpZiBuZXR3b3JrLmlzX29ubGluZSgpOiBzbGVlcCgxKSIsCiAgICAgICAgICAgICJ3aGlsZSBuZXR3b3JrLmlzX29ubGluZSgpOiBzbGVlcCgxKSIsCiAgICAgICAgICAgICJ3aGlsZSBub3QgbmV0d29yay5pc19vbmxpbmUoKTogc2xlZXAoMSkiLAogICAgICAgICAgICAiaWYgbm90IG5ldHdvcmsuaXNfb25saW5lKCk6IHNsZWVwKDEpIgogICAgICAgIF0sCiAgICAgICAgInNjb3JlIjogIndoaWxlIG5vdCBuZXR3b3JrLmlzX29ubGluZSgpOiBzbGVlcCgxKSIKICAgIH0sCiAgICB7CiAgICAgICAgInRpdGxlIjogIllvdSBhcmUgY29sbGFib3JhdGluZyBvbiBhIGNvZGViYXNlIHRoYXQgdXNlcyBxdWl0ZSBhIGZldyBmdW5jdGlvbnMgdGhhdCBhcmUgY29waWVkIGluIGEgZmV3IG90aGVyIGZpbGVzLCBzbyB5b3UgZGVjaWRlIHRvIHJlZmFjdG9yIHRoZW0gaW50byBhIGxpYnJhcnkuIFRoZSBsZWFkIGRldiB3YW50cyB5b3UgdG8gcHJvdmlkZSB5b3VyIHJhdGlvbmFsZSBmb3IgdGhpcyByZWZhY3Rvcml6YXRpb24gaW4gdGhlIHB1bGwgcmVxdWVzdC5XaGF0IGJlbmVmaXRzIHNob3VsZCB5b3UgcG9pbnQgb3V0IGZvciBtb3ZpbmcgdGhlIGNvZGUgaW50byBhIGxpYnJhcnk/IiwKICAgICAgICAicXVlc190eXBlIjogMTUsCiAgICAgICAgIm9wdGlvbnMiOiBbCiAgICAgICAgICAgICJJbmNyZWFzZSBwcml2YWN5LiIsCiAgICAgICAgICAgICJSZWR1Y2UgY29tcGlsZXIgZXJyb3JzLiIsCiAgICAgICAgICAgICJJbXByb3ZlIHJ1bnRpbWUgc3BlZWQuIiwKICAgICAgICAgICAgIkNyZWF0ZSBhIG5hbWVzcGFjZSBmb3IgbGlicmFyeSBmdW5jdGlvbnMuIiwKICAgICAgICAgICAgIk5hdmlnYXRlIGNvZGViYXNlIG1vcmUgZWFzaWx5LiIKICAgICAgICBdLAogICAgICAgICJzY29yZSI6IFsKICAgICAgICAgICAgIkNyZWF0ZSBhIG5hbWVzcGFjZSBmb3IgbGlicmFyeSBmdW5jdGlvbnMuIiwKICAgICAgICAgICAgIk5hdmlnYXRlIGNvZGViYXNlIG1vcmUgZWFzaWx5LiIKICAgICAgICBdCiAgICB9LAogICAgewogICAgICAgICJ0aXRsZSI6ICJZb3UgYXJlIHdvcmtpbmcgb24gYSBncmFwaGljYWwgaW50ZXJmYWNlIGZvciBhIGZpbGUgZXhwbG9yZXIuIFRoZSBmaWxlIGV4cGxvcmVyXHUyMDE5cyBncmFwaGljYWwgd2luZG93IGlzIG1hbmFnZWQgYnkgdHdvIGNsYXNzZXM6IG9uZSB0byBtYW5hZ2UgY2hhbmdlcyBpbiB0aGUgd2luZG93LCBhbmQgYW5vdGhlciB0byBkZWFsIHdpdGggaW5wdXQgZXZlbnRzLiBEdXJpbmcgcmV2aWV3LCB5b3VyIG1hbmFnZXIgYXJndWVzIHRoYXQgdGhpcyBkZXNpZ24gdmlvbGF0ZXMgdGhlIFNSUCAoc2luZ2xlIHJlc3BvbnNpYmlsaXR5IHByaW5jaXBsZSkuSG93IGNhbiB5b3UgcmVmYWN0b3IgdGhlIGNvZGUgdG8gcmVzcGVjdCB0aGUgU1JQPyIsCiAgICAgICAgInF1ZXNfdHlwZSI6IDIsCiAgICAgICAgIm9wdGlvbnMiOiBbCiAgICAgICAgICAgICJNb3ZlIGVhY2ggY2xhc3MgdG8gaXRzIG93biBwcm9qZWN0IGZpbGUuIiwKICAgICAgICAgICAgIkNyZWF0ZSBhbiBhZGRpdGlvbmFsIGNsYXNzIHRvIHBhcnNlIGV2ZW50cyBhbmQgcGFzcyB0aGVtIHRvIHRoZSB3aW5kb3cuIiwKICAgICAgICAgICAgIkRlY29uc3RydWN0IHRoZSBjbGFzcyBpbnRvIHNpbXBsZSBmdW5jdGlvbnMuIiwKICAgICAgICAgICAgIkNyZWF0ZSBhIHNpbmdsZSBjbGFzcyB0byBtYW5hZ2UgdGhlIGZpbGUgZXhwbG9yZXIuIgogICAgICAgIF0sCiAgICAgICAgInNjb3JlIjogIkNyZWF0ZSBhIHNpbmdsZSBjbGFzcyB0byBtYW5hZ2UgdGhlIGZpbGUgZXhwbG9yZXIuIgogICAgfQpd"
$bytes = [System.Convert]::FromBase64String($b64)
$newText = [System.Text.Encoding]::UTF8.GetString($bytes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the shared-string payload before; remove it (and its row) entirely.
$ws.Range("A2").ClearContents()

# A1 held a plain numeric 0 with a bold/bordered/centered style; put the
# payload text there instead and reset formatting back to the workbook
# default ("Normal") so no custom font/border survives on the cell.
$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"
